# Any file with Ra_Stock_5 used had the wrong concentration for stock 5
# radium; correct the Stock Activity (value + error) on the Parameters sheet.
# Was B6=675.45997631774435 / C6=0.13742948708443162 -> B6=1407 / C6=62.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Activate()

$ws.Range("B6").Value = 1407
$ws.Range("C6").Value = 62
$ws.Range("C6").NumberFormat = "0.00E+00"

$ws.Range("B6:C6").Select()

# The calibration scatter chart's Activity (Y) axis is displayed in
# scientific notation now that the corrected stock activity flows through.
$calSheet = $wb.Worksheets.Item("Calibration Data")
$chart = $calSheet.ChartObjects().Item(1).Chart
$chart.Axes(2).NumberFormat = "0.00E+00"
